$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions) - simple numeric / text refreshes
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1499
$ws1.Range("F6").Value  = 775
$ws1.Range("F7").Value  = 50
$ws1.Range("F12").Value = 36963
$ws1.Range("G12").Value = "暂时售罄"
$ws1.Range("F13").Value = 7420
$ws1.Range("F15").Value = 398
$ws1.Range("F16").Value = 600
$ws1.Range("F19").Value = 126
$ws1.Range("F20").Value = 460
$ws1.Range("F23").Value = 475
$ws1.Range("F24").Value = 140
$ws1.Range("F25").Value = 857
$ws1.Range("F29").Value = 470
$ws1.Range("F32").Value = 74
$ws1.Range("F36").Value = 80
$ws1.Range("F37").Value = 783

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1251

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1526
$ws3.Range("F3").Value = 378

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - combined listing
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

# simple refreshes that don't involve the row-14 insert / row-20 removal
$ws4.Range("F2").Value  = 1526
$ws4.Range("F3").Value  = 378
$ws4.Range("F4").Value  = 1251
$ws4.Range("F7").Value  = 1499
$ws4.Range("F9").Value  = 775
$ws4.Range("F10").Value = 50

# A new exhibition ("萤火虫动漫游戏嘉年华 x KKWORLD2024") was published and now
# needs to be listed at row 14; every later row (14-20) shifts down by one,
# and the old row 20 ("AP动漫游戏嘉年华·徐慧内场") drops off the list entirely.
# Column A is a plain positional index (row-1) and is rewritten explicitly so
# it keeps matching its row after the contents below are shuffled.

$row14 = @("2024-07-19", "广州·萤火虫动漫游戏嘉年华 × KKWORLD2024 快看漫画乐园", "新港东路1000号 保利世贸博览馆", "2024.07.19 09:00-07.22 17:00", 36963, "暂时售罄", "https://show.bilibili.com/platform/detail.html?id=87210", "//i1.hdslb.com/bfs/openplatform/202406/DTCdOTPs1718177177472.jpeg")
$row15 = @("2024-07-20", "广州·冰兔2024线下live「过去和未来」", "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）", "2024.07.20 20:00-07.20 22:00", 175, 198, "https://show.bilibili.com/platform/detail.html?id=87546", "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg")
$row16 = @("2024-07-20", "广州·跨越二次元ACG神级动漫世界巡回演唱会", "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院", "2024.07.20 19:30-07.20 21:10", 301, 280, "https://show.bilibili.com/platform/detail.html?id=85353", "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg")
$row17 = @("2024-07-21", "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024", "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse", "2024.07.21 14:30-07.21 16:00", 249, 280, "https://show.bilibili.com/platform/detail.html?id=87034", "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png")
$row18 = @("2024-07-26", "广州·【早鸟8折】“浪漫古典Ⅱ”百年经典传世名曲烛光音乐会 ", "广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）", "2024.07.26 20:00-07.26 21:30", 5, 144, "https://show.bilibili.com/platform/detail.html?id=87726", "//i0.hdslb.com/bfs/openplatform/202406/A8vhVlhn1717575084179.png")
$row19 = @("2024-07-26", "广州·萨克斯王子安德鲁·杨——2024经典&流行音乐巡回演出", "龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House", "2024.07.26 20:00-07.26 21:30", 5, 380, "https://show.bilibili.com/platform/detail.html?id=86635", "//i1.hdslb.com/bfs/openplatform/202405/rciNih361716802006584.jpeg")
$row20 = @("2024-07-27", "广州·AP动漫游戏嘉年华", "新港东路630-638号 南丰国际会展中心", "2024.07.27 09:00-07.28 17:00", 7420, 80, "https://show.bilibili.com/platform/detail.html?id=87213", "//i1.hdslb.com/bfs/openplatform/202406/3Z8rGZPP1718164976101.jpeg")

$newRows = @($row14, $row15, $row16, $row17, $row18, $row19, $row20)

# Column B holds plain "yyyy-mm-dd" text, not real dates - force Text format
# first so Excel doesn't silently reinterpret the strings as date serials.
$ws4.Range("B14:B20").NumberFormat = "@"

$r = 14
foreach ($data in $newRows) {
    $ws4.Range("A$r").Value = $r - 1
    $ws4.Range("B$r").Value = $data[0]
    $ws4.Range("C$r").Value = $data[1]
    $ws4.Range("D$r").Value = $data[2]
    $ws4.Range("E$r").Value = $data[3]
    $ws4.Range("F$r").Value = $data[4]
    $ws4.Range("G$r").Value = $data[5]
    $ws4.Range("H$r").Value = $data[6]
    $ws4.Range("I$r").Value = $data[7]
    $r = $r + 1
}

# row 21 (LookLook动漫嘉年华2th) keeps its identity, just the view counter ticks up
$ws4.Range("F21").Value = 398

# remaining simple refreshes below the shifted block
$ws4.Range("F23").Value = 600
$ws4.Range("F27").Value = 126
$ws4.Range("F28").Value = 460
$ws4.Range("F31").Value = 475
$ws4.Range("F32").Value = 140
$ws4.Range("F33").Value = 857
$ws4.Range("F37").Value = 470
$ws4.Range("F40").Value = 74
